# The deck's theme parts were swapped: the "Office Theme" palette (formerly
# ppt/theme/theme1.xml, used only by the notes master) becomes the palette
# of the presentation's active theme (ppt/theme/theme2.xml, used by the
# slide master / ActivePresentation), replacing the "Integral" palette that
# lived there. Fonts and format scheme are identical between the two themes,
# so only the twelve theme colors need to change.
#
# PowerPoint's ColorScheme object is the COM surface for editing a theme's
# clrScheme in place: Colors(index).RGB, indexed 1-12 in clrScheme document
# order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). RGB values use the
# standard VBA encoding R + G*256 + B*65536.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme

$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
